# Apply cryptos list update (prices + 1h volume %) for rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.091.11"
$ws.Range("E2").Value = "  -3.92%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.496.14"
$ws.Range("E3").Value = "  -5.95%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.45"
$ws.Range("E5").Value = "  -1.51%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.67"
$ws.Range("E6").Value = "  -4.67%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.615"
$ws.Range("E7").Value = "  -0.35%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.487.89"
$ws.Range("E8").Value = "  -5.99%  "

$ws.Range("E9").Value = "  +0.13%  "

$ws.Range("E10").Value = "  -7.81%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.67"
$ws.Range("E11").Value = "  +4.07%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.592"
$ws.Range("E12").Value = "  -3.59%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.56"
$ws.Range("E13").Value = "  -6.96%  "

$ws.Range("E14").Value = "  -4.99%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "672.76"
$ws.Range("E15").Value = "  -1.47%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.060.11"
$ws.Range("E16").Value = "  -5.96%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.66"
$ws.Range("E17").Value = "  -4.33%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.098.78"
$ws.Range("E18").Value = "  -4.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.502.22"
$ws.Range("E19").Value = "  -5.48%  "

$ws.Range("E20").Value = "  -1.57%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.35"
$ws.Range("E21").Value = "  -4.49%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.12"
$ws.Range("E22").Value = "  -5.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.897"
$ws.Range("E23").Value = "  -5.44%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.04"
$ws.Range("E24").Value = "  -10.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.28"
$ws.Range("E25").Value = "  -6.38%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.85"
$ws.Range("E26").Value = "  -5.11%  "

$ws.Range("E27").Value = "  -0.06%  "

$ws.Range("E28").Value = "  -7.47%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.36"
$ws.Range("E29").Value = "  -9.83%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.79"
$ws.Range("E30").Value = "  -8.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.64"
$ws.Range("E31").Value = "  -7.40%  "

$ws.Range("E32").Value = "  -8.84%  "

$ws.Range("E33").Value = "  -6.39%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.21"
$ws.Range("E34").Value = "  -2.34%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "595.05"
$ws.Range("E35").Value = "  +5.21%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.57"
$ws.Range("E36").Value = "  -15.23%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.83"
$ws.Range("E37").Value = "  -4.28%  "

$ws.Range("E38").Value = "  -5.65%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "56.91"
$ws.Range("E39").Value = "  -4.72%  "

$ws.Range("E40").Value = "  +0.11%  "

$ws.Range("E41").Value = "  -6.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.333"
$ws.Range("E42").Value = "  -5.73%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.405.86"
$ws.Range("E43").Value = "  -9.74%  "

$ws.Range("E44").Value = "  -6.69%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "33.12"
$ws.Range("E45").Value = "  -7.73%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₃0703"
$ws.Range("E46").Value = "  -10.03%  "

$ws.Range("E47").Value = "  -1.00%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.58"
$ws.Range("E48").Value = "  -8.40%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.132"
$ws.Range("E49").Value = "  -1.25%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.74"
$ws.Range("E50").Value = "  +16.64%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "132.90"
$ws.Range("E51").Value = "  -2.22%  "
